$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 (Beta) values ---
$ws.Range("C2").Value = 16.27395536621195
$ws.Range("E2").Value = 0.006275575909959944
$ws.Range("F2").Value = 6.608387435978644
$ws.Range("G2").Value = 6.398070021960509
$ws.Range("H2").Value = 6.823114605078502
$ws.Range("I2").Value = 3.2994696946754
$ws.Range("J2").Value = 3.270501540739864
$ws.Range("K2").Value = 3.330710973808486
$ws.Range("L2").Value = 0.1836414946648398
$ws.Range("M2").Value = 0.1820573890941121
$ws.Range("N2").Value = 0.1853512161757067

# --- Update row 3 (Gamma) values ---
$ws.Range("C3").Value = 0.3589057182506037
$ws.Range("D3").Value = 0.3038166771491592
$ws.Range("E3").Value = 0.3557873748505794
$ws.Range("F3").Value = 0.2594165550832597
$ws.Range("G3").Value = 0.2590128726673701
$ws.Range("H3").Value = 0.2598164956945308
$ws.Range("I3").Value = 0.2326383659222323
$ws.Range("J3").Value = 0.2322669381939036
$ws.Range("K3").Value = 0.2330025539085543
$ws.Range("L3").Value = 0.2569356968891153
$ws.Range("M3").Value = 0.2565347254353867
$ws.Range("N3").Value = 0.2573330564668757

# --- Add new row 4 (Beta + Gamma) ---
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Beta + Gamma"
$ws.Range("C4").Value = 16.63286108446255
$ws.Range("D4").Value = 0.3068349094013313
$ws.Range("E4").Value = 0.3620629507605393
$ws.Range("F4").Value = 6.867803991061905
$ws.Range("G4").Value = 6.657082894627879
$ws.Range("H4").Value = 7.082931100773033
$ws.Range("I4").Value = 3.532108060597632
$ws.Range("J4").Value = 3.502768478933767
$ws.Range("K4").Value = 3.56371352771704
$ws.Range("L4").Value = 0.4405771915539552
$ws.Range("M4").Value = 0.4385921145294988
$ws.Range("N4").Value = 0.4426842726425823
